$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.166.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.899.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3908"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07878"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9893"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.883.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.38%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.060"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.82%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.731"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07003"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.79%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.70%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.08%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009976"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.32%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.186.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.302"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.129.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.54%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.110"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.94%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.71%  "

# Row 27
$ws.Range("E27").Value = "  -0.47%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.946"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.45%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.49%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.879"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09340"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8994"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.87%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.246"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.90%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.323"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.96%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.176"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.70%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.179"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05760"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02084"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.00%  "

# Row 39
$ws.Range("E39").Value = "  -0.15%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.716"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.53%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5692"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.07%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1788"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.701"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.69%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5344"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.168"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.30%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07011"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.851"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.41%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.555"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.34%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.60%  "

# Row 51
$ws.Range("E51").Value = "  -1.24%  "
